$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead text: volume/number and week-of date range ---
$ws.Range("A8").Value = "Volume 32   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/13/2025  Through  10/19/2025"

# --- Column widths for I (9) and J (10): widen to match new wider numbers ---
$ws.Columns.Item(9).ColumnWidth = 7.433768
$ws.Columns.Item(10).ColumnWidth = 7.433768

# --- Re-type cells whose underlying value flips between a number and the
#     "0"/"***.*" placeholder text (copy a same-styled neighbour first so the
#     destination inherits the right style + data type, then stamp the value). ---
$ws.Range("C29").Copy($ws.Range("C14"))
$ws.Range("C29").Copy($ws.Range("C15"))
$ws.Range("D29").Copy($ws.Range("D22"))
$ws.Range("E29").Copy($ws.Range("E22"))
$ws.Range("C29").Copy($ws.Range("C27"))
$ws.Range("D16").Copy($ws.Range("D28"))
$ws.Range("E16").Copy($ws.Range("E28"))
$ws.Range("D16").Copy($ws.Range("D31"))
$ws.Range("E16").Copy($ws.Range("E31"))
$ws.Range("G16").Copy($ws.Range("G31"))
$ws.Range("H16").Copy($ws.Range("H31"))
$ws.Range("D16").Copy($ws.Range("D33"))
$ws.Range("E16").Copy($ws.Range("E33"))
$ws.Range("G16").Copy($ws.Range("G33"))
$ws.Range("H16").Copy($ws.Range("H33"))

# --- Final cell values for rows 14-33 ---
# Row 14
$ws.Range("N14").Value = -66.666666666666
# Row 15
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 31
$ws.Range("J15").Value = 34
$ws.Range("K15").Value = -8.823529411764
$ws.Range("L15").Value = 14.814814814814
$ws.Range("M15").Value = 138.461538461538
$ws.Range("N15").Value = -53.030303030303
# Row 16
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -42.857142857142
$ws.Range("F16").Value = 32
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = 6.666666666666
$ws.Range("I16").Value = 364
$ws.Range("J16").Value = 367
$ws.Range("K16").Value = -0.817438692098
$ws.Range("L16").Value = 5.813953488372
$ws.Range("M16").Value = 48.571428571428
$ws.Range("N16").Value = -63.960396039604
# Row 17
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 40
$ws.Range("F17").Value = 53
$ws.Range("G17").Value = 51
$ws.Range("H17").Value = 3.921568627450
$ws.Range("I17").Value = 634
$ws.Range("J17").Value = 602
$ws.Range("K17").Value = 5.315614617940
$ws.Range("L17").Value = 11.228070175438
$ws.Range("M17").Value = 145.736434108527
$ws.Range("N17").Value = -20.050441361916
# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = -44.444444444444
$ws.Range("I18").Value = 188
$ws.Range("J18").Value = 234
$ws.Range("K18").Value = -19.658119658119
$ws.Range("L18").Value = 13.939393939393
$ws.Range("M18").Value = 49.206349206349
$ws.Range("N18").Value = -78.636363636363
# Row 19
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 22.580645161290
$ws.Range("I19").Value = 474
$ws.Range("J19").Value = 440
$ws.Range("K19").Value = 7.727272727272
$ws.Range("L19").Value = 13.126491646778
$ws.Range("M19").Value = 143.076923076923
$ws.Range("N19").Value = 59.595959595959
# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = 28.571428571428
$ws.Range("I20").Value = 213
$ws.Range("J20").Value = 225
$ws.Range("K20").Value = -5.333333333333
$ws.Range("L20").Value = -40.336134453781
$ws.Range("M20").Value = 131.521739130435
$ws.Range("N20").Value = -43.650793650793
# Row 21
$ws.Range("C21").Value = 35
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -7.894736842105
$ws.Range("F21").Value = 168
$ws.Range("G21").Value = 162
$ws.Range("H21").Value = 3.703703703703
$ws.Range("I21").Value = 1913
$ws.Range("J21").Value = 1907
$ws.Range("K21").Value = 0.314630309386
$ws.Range("L21").Value = 1.109936575052
$ws.Range("M21").Value = 103.944562899787
$ws.Range("N21").Value = -44.566792234135
# Row 22
$ws.Range("L22").Value = -66.666666666666
# Row 23
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 19
$ws.Range("G23").Value = 21
$ws.Range("H23").Value = -9.523809523809
$ws.Range("I23").Value = 242
$ws.Range("J23").Value = 337
$ws.Range("K23").Value = -28.189910979228
$ws.Range("L23").Value = -28.189910979228
$ws.Range("M23").Value = 58.169934640522
# Row 24
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -20.833333333333
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = -8.988764044943
$ws.Range("I24").Value = 949
$ws.Range("J24").Value = 831
$ws.Range("K24").Value = 14.199759326113
$ws.Range("L24").Value = -0.315126050420
$ws.Range("M24").Value = 52.572347266881
# Row 25
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 0
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = -46.666666666666
$ws.Range("I25").Value = 170
$ws.Range("J25").Value = 142
$ws.Range("K25").Value = 19.718309859154
$ws.Range("L25").Value = -21.658986175115
# Row 26
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 13
$ws.Range("F26").Value = 66
$ws.Range("G26").Value = 85
$ws.Range("H26").Value = -22.352941176470
$ws.Range("I26").Value = 731
$ws.Range("J26").Value = 792
$ws.Range("K26").Value = -7.702020202020
$ws.Range("L26").Value = -15
$ws.Range("M26").Value = 2.957746478873
# Row 27
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 36
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = -28
$ws.Range("L27").Value = -21.739130434782
# Row 28
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 25
$ws.Range("I28").Value = 48
$ws.Range("J28").Value = 66
$ws.Range("K28").Value = -27.272727272727
$ws.Range("L28").Value = -40.740740740740
# Row 29
$ws.Range("N29").Value = -68.478260869565
# Row 30
$ws.Range("N30").Value = -70.329670329670
# Row 31
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = -100
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = -100
$ws.Range("J31").Value = 2
$ws.Range("K31").Value = 0
# Row 33
$ws.Range("D33").Value = 2
$ws.Range("E33").Value = -100
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = -100
$ws.Range("J33").Value = 5
$ws.Range("K33").Value = -60
